$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5
$ws.Range("A5").Value = 111482955
$ws.Range("I5").Value = '70'
$ws.Range("P5").Value = 'S om järnvägen - 5, Vg'
$ws.Range("Q5").Value = 432064.1298546481
$ws.Range("R5").Value = 6419677.395781181

# Row 6
$ws.Range("A6").Value = 111482936
$ws.Range("I6").Value = '25'
$ws.Range("P6").Value = 'S om järnvägen - 4, Vg'
$ws.Range("Q6").Value = 432073.5656663703
$ws.Range("R6").Value = 6419668.734013095

# Row 8
$ws.Range("A8").Value = 111483037
$ws.Range("B8").Value = 96348
$ws.Range("D8").Value = 'VU'
$ws.Range("E8").Value = 220787
$ws.Range("F8").Value = 'Knärot'
$ws.Range("G8").Value = 'Goodyera repens'
$ws.Range("H8").Value = '(L.) R. Br.'
$ws.Range("I8").Value = '60'
$ws.Range("J8").Value = 'stjälkar/strån/skott'
$ws.Range("K8").Value = 'blomning'
$ws.Range("L8").Value = ''
$ws.Range("P8").Value = 'S om järnvägen - 7, Vg'
$ws.Range("Q8").Value = 432060.6482816387
$ws.Range("R8").Value = 6419660.45125766
$ws.Range("AJ8").Value = ""
$ws.Range("AK8").Value = ""
$ws.Range("AM8").Value = ""
$ws.Range("AO8").Value = ""

# Row 9
$ws.Range("A9").Value = 111483381
$ws.Range("P9").Value = 'S om järnvägen - 14, Vg'
$ws.Range("Q9").Value = 431754.10213514
$ws.Range("R9").Value = 6419728.893211351

# Row 10
$ws.Range("A10").Value = 111483105
$ws.Range("B10").Value = 73689
$ws.Range("D10").Value = 'NT'
$ws.Range("E10").Value = 308
$ws.Range("F10").Value = 'Brunpudrad nållav'
$ws.Range("G10").Value = 'Chaenotheca gracillima'
$ws.Range("H10").Value = '(Vain.) Tibell'
$ws.Range("I10").Value = ''
$ws.Range("J10").Value = ''
$ws.Range("K10").Value = ''
$ws.Range("L10").Value = ""
$ws.Range("P10").Value = 'S om järnvägen - 8, Vg'
$ws.Range("Q10").Value = 431947.1499479365
$ws.Range("R10").Value = 6419623.056550305
$ws.Range("AJ10").Value = 'tall'
$ws.Range("AK10").Value = 'Pinus sylvestris'
$ws.Range("AM10").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO10").Value = 'Standing dead tree/snags # Pinus sylvestris'

# Row 11
$ws.Range("A11").Value = 111491635
$ws.Range("I11").Value = '10'
$ws.Range("K11").Value = 'blomning'
$ws.Range("P11").Value = 'S om järnvägen - 21, Vg'
$ws.Range("Q11").Value = 431859.6228004749
$ws.Range("R11").Value = 6419672.898494411

# Row 12
$ws.Range("A12").Value = 111490843
$ws.Range("B12").Value = 96348
$ws.Range("D12").Value = 'VU'
$ws.Range("E12").Value = 220787
$ws.Range("F12").Value = 'Knärot'
$ws.Range("G12").Value = 'Goodyera repens'
$ws.Range("H12").Value = '(L.) R. Br.'
$ws.Range("I12").Value = '50'
$ws.Range("J12").Value = 'stjälkar/strån/skott'
$ws.Range("K12").Value = 'fullt utvecklade blad'
$ws.Range("L12").Value = ''
$ws.Range("P12").Value = 'S om järnvägen - 17, Vg'
$ws.Range("Q12").Value = 431803.2980747336
$ws.Range("R12").Value = 6419679.170503675
$ws.Range("AM12").Value = ""
$ws.Range("AO12").Value = ""

# Row 13
$ws.Range("A13").Value = 111483437
$ws.Range("I13").Value = '100'
$ws.Range("P13").Value = 'S om järnvägen - 15, Vg'
$ws.Range("Q13").Value = 431797.479853621
$ws.Range("R13").Value = 6419681.394993878

# Row 14
$ws.Range("A14").Value = 111483197
$ws.Range("B14").Value = 73689
$ws.Range("D14").Value = 'NT'
$ws.Range("E14").Value = 308
$ws.Range("F14").Value = 'Brunpudrad nållav'
$ws.Range("G14").Value = 'Chaenotheca gracillima'
$ws.Range("H14").Value = '(Vain.) Tibell'
$ws.Range("I14").Value = ''
$ws.Range("J14").Value = ''
$ws.Range("K14").Value = ''
$ws.Range("L14").Value = ""
$ws.Range("P14").Value = 'S om järnvägen - 11, Vg'
$ws.Range("Q14").Value = 431937.082796899
$ws.Range("R14").Value = 6419625.884406033
$ws.Range("AM14").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO14").Value = 'Standing dead tree/snags'

# Row 15
$ws.Range("A15").Value = 111483300
$ws.Range("B15").Value = 73689
$ws.Range("D15").Value = 'NT'
$ws.Range("E15").Value = 308
$ws.Range("F15").Value = 'Brunpudrad nållav'
$ws.Range("G15").Value = 'Chaenotheca gracillima'
$ws.Range("H15").Value = '(Vain.) Tibell'
$ws.Range("I15").Value = ''
$ws.Range("J15").Value = ''
$ws.Range("K15").Value = ''
$ws.Range("L15").Value = ""
$ws.Range("P15").Value = 'S om järnvägen - 12, Vg'
$ws.Range("Q15").Value = 431888.091041417
$ws.Range("R15").Value = 6419625.122914318
$ws.Range("AJ15").Value = 'tall'
$ws.Range("AK15").Value = 'Pinus sylvestris'
$ws.Range("AM15").Value = 'Stående död trädstam/högstubbe'
$ws.Range("AO15").Value = 'Standing dead tree/snags # Pinus sylvestris'

# Row 16
$ws.Range("A16").Value = 111483462
$ws.Range("B16").Value = 96348
$ws.Range("D16").Value = 'VU'
$ws.Range("E16").Value = 220787
$ws.Range("F16").Value = 'Knärot'
$ws.Range("G16").Value = 'Goodyera repens'
$ws.Range("H16").Value = '(L.) R. Br.'
$ws.Range("I16").Value = '45'
$ws.Range("J16").Value = 'stjälkar/strån/skott'
$ws.Range("K16").Value = 'blomning'
$ws.Range("L16").Value = ''
$ws.Range("P16").Value = 'S om järnvägen - 16, Vg'
$ws.Range("Q16").Value = 431654.0242198514
$ws.Range("R16").Value = 6419791.70470859
$ws.Range("AJ16").Value = ""
$ws.Range("AK16").Value = ""
$ws.Range("AM16").Value = ""
$ws.Range("AO16").Value = ""

# Row 17
$ws.Range("A17").Value = 111491187
$ws.Range("I17").Value = '60'
$ws.Range("P17").Value = 'S om järnvägen - 18, Vg'
$ws.Range("Q17").Value = 431829.514510141
$ws.Range("R17").Value = 6419749.394753682

# Row 18
$ws.Range("A18").Value = 111483107
$ws.Range("B18").Value = 73681
$ws.Range("D18").Value = 'LC'
$ws.Range("E18").Value = 6439
$ws.Range("F18").Value = 'Gulnål'
$ws.Range("G18").Value = 'Chaenotheca brachypoda'
$ws.Range("H18").Value = '(Ach.) Tibell'
$ws.Range("P18").Value = 'S om järnvägen - 8, Vg'
$ws.Range("Q18").Value = 431947.1499479365
$ws.Range("R18").Value = 6419623.056550305
